# Remplacement des montagnes par des hexagones retirés
$wb = $excel.ActiveWorkbook

# Rename the "Montagnes" sheet to "Obstacles"
$ws = $wb.Worksheets.Item("Montagnes")
$ws.Name = "Obstacles"

# The header cell on this sheet mirrors the sheet's subject -> update its text too
$ws.Range("A2").Value = "Obstacles"

# Make "Obstacles" the active sheet/tab and set its selection
$ws.Activate()
$ws.Range("E19").Select()

$wb.Save()
